$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header D1: "Tempo de uso" -> "Tempo de uso em segundos"
$ws.Range("D1").Value = "Tempo de uso em segundos"

# Column D needs to be widened to fit the new, longer header text
$ws.Columns(4).ColumnWidth = 29

# E4 had a stray duplicate number-format style; normalize it to the same
# time format (hh:mm:ss) used by the other cells in column E
$ws.Range("E4").NumberFormat = "hh:mm:ss"

# Update the saved cursor/selection position
$ws.Range("E11").Select()
